$wb = $excel.ActiveWorkbook

# "soil" sheet: area_exceeded_km2 values recomputed with a wider (32-bit)
# counter, fixing the 16-bit overflow in the upstream pipeline.
$soil = $wb.Worksheets.Item("soil")
$soil.Range("C2").Value = 1379
$soil.Range("C3").Value = 1848
$soil.Range("C4").Value = 796
$soil.Range("C5").Value = 199
$soil.Range("C6").Value = 0
$soil.Range("C7").Value = 0
$soil.Range("C8").Value = 0
$soil.Range("C9").Value = 0
$soil.Range("C10").Value = 0
$soil.Range("C11").Value = 0

# Update which sheet/cell is active & selected: "vegetation" loses focus,
# "soil" becomes the active tab with C12 selected.
$veg = $wb.Worksheets.Item("vegetation")
$veg.Activate()
$veg.Range("D19").Select()

$soil.Activate()
$soil.Range("C12").Select()
